$d = $word.ActiveDocument

# --- 1) Title paragraph + following empty paragraph: drop the nl-BE language overrides ---
$pTitle = $d.Paragraphs.Item(1)
$pAfterTitle = $d.Paragraphs.Item(2)
$rTitle = $d.Range($pTitle.Range.Start, $pAfterTitle.Range.End)
$titleXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>ReadMe</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
"@
$null = $rTitle.InsertXML($titleXml)

# --- 2) Final paragraph: rewrite the sentence and split off the bookmark into its own paragraph ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$lastXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>In this repo you will find al</w:t></w:r><w:r><w:t>l the supplementary material that is linked to</w:t></w:r><w:r><w:t xml:space="preserve"> my master thesis.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$null = $rLast.InsertXML($lastXml)
